$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "AdServingType" column (J) with header and a few data values
$ws.Range("J1").Value = "AdServingType"
$ws.Range("J2").Value = "Click Tracker"
$ws.Range("J6").Value = "Tracking"
$ws.Range("J7").Value = "No Tracking"

# Select J8 to match the final selection state
$ws.Range("J8").Select()
